$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

$ws.Range("O2").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Aceptaci%E0%B8%82n/Aceptaciขn+1+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Aceptaci%E0%B8%82n/Aceptaciขn+2+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Aceptaci%E0%B8%82n/Aceptaciขn+3+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Aceptaci%E0%B8%82n/Aceptaciขn+4+Rojo.jpg']"

$ws.Range("O3").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Afirmaci%E0%B8%82n/Afirmaciขn+1+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Afirmaci%E0%B8%82n/Afirmaciขn+2+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Afirmaci%E0%B8%82n/Afirmaciขn+3+Negro.jpg']"

$ws.Range("O8").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+1+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+2+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+3+Negro.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+4+Azul.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+5+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Somos+La+Verdad/Intuici%E0%B8%82n/Intuiciขn+6+Rojo.png']"
